$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.238.07"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.603.15"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.66"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3778"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.76"
$ws.Range("E8").Value = "  +3.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3641"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.275"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08147"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.84"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.426"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001250"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "1.606.86"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.10"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06925"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.538"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "23.246.51"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.027"
$ws.Range("E25").Value = "  +8.99%  "
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.25"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.51"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.262"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.37"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.771"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").Value = "1.782.83"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9653"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07542"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02748"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.26"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2541"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.131"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08827"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.389"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7119"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.56"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.70"
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6548"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.327"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9993"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.016"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.77"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07957"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.207"
$ws.Range("E51").Value = "  -2.42%  "
